# Fix relationship for a couple of categories (euklems-2018-isic4-to-fingreen map)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 31-33 (ISIC "31-33" -> fingreen C31_32 / C33) were mis-tagged as
# "recomposition"; they are actually a one-to-many split -> "disaggregation".
$ws.Range("C23").Value = "disaggregation"
$ws.Range("C24").Value = "disaggregation"

# Drop stray formatting-only (empty) cells left in column B for the "G" and
# "J" aggregation header rows.
$ws.Range("B30").ClearContents()
$ws.Range("B39").ClearContents()

# Fix the "recomposition_method" for the S row (inherit), and add a new row
# so both S and T recompose (average) into the combined "ST" fingreen code,
# in addition to S's existing "S95" mapping.
$ws.Range("D54").Value = "inherit"

$ws.Rows("55").Insert()
$ws.Range("A55").Value = "S"
$ws.Range("B55").Value = "ST"
$ws.Range("C55").Value = "recomposition"
$ws.Range("D55").Value = "average"

# The old row 55 (T / ST / recomposition) shifted down to row 56; it was
# missing its recomposition_method -> "average".
$ws.Range("D56").Value = "average"

[void]$ws.Range("C24").Select()
